$d = $word.ActiveDocument

# Locate the end of the final paragraph's text (right where the
# existing _GoBack bookmark sits, immediately before the paragraph
# mark). We need to split this paragraph in two and insert a new
# paragraph of text there, with the bookmark ending up at the end of
# that new paragraph.
$bm = $d.Bookmarks("_GoBack")
$splitPos = $bm.Range.Start
$null = $bm.Delete()

$insertionPoint = $d.Range($splitPos, $splitPos)

$newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Ask </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Venkata</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> if the current behavior for the genes of interest is okay. Right now, selecting MCM3-S, gives us DMNT1-E as one of the children. DMNT1-E has a degree of 17, so since we are showing top 3 interactions, one would expect that there would be 4 edges involving the node for DMNT1-E: 1 for the edge coming from MCM3-S, and 3 for the next top 3 interactions. However, there are currently only 3 edges showing since MCM3-S is one of the top 3 interactors. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $insertionPoint.InsertXML($newParagraphXml)

# The document previously ended with a trailing empty paragraph; that
# is now redundant since the new paragraph above takes its place, so
# remove it by deleting from the end of the new last paragraph's text
# through the end of the document (merges away the stray empty
# paragraph mark).
$lastContentPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$tailRange = $d.Range($lastContentPara.Range.End - 1, $d.Content.End)
$tailRange.Delete()
